# Update values on Sheet1 to reflect re-pulled data / pushed all data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = -1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 9

# Row 8
$ws.Range("F8").Value = -3

# Row 10
$ws.Range("F10").Value = -1

# Row 16
$ws.Range("F16").Value = 0

# Row 22
$ws.Range("F22").Value = 1

# Row 39
$ws.Range("F39").Value = -2

# Row 45
$ws.Range("F45").Value = 2

# Row 50
$ws.Range("F50").Value = 3

# Row 51
$ws.Range("F51").Value = 2

# Row 53
$ws.Range("F53").Value = 0

# Row 57
$ws.Range("F57").Value = 0

# Row 63
$ws.Range("F63").Value = 1

# Row 68
$ws.Range("F68").Value = -2

# Row 69
$ws.Range("F69").Value = -3

# Row 73
$ws.Range("F73").Value = -3
